# Fix contact information missing from short resumes.
# Insert a new centered paragraph with contact info directly after the
# "Dheeraj Chand" name heading, matching the long-resume layout.
#
# Using Find/Replace with a paragraph-mark wildcard (^p) rather than
# Range.InsertParagraphAfter() keeps the new paragraph's run free of the
# bold/28pt character formatting that the name run carries (InsertParagraphAfter
# copies the formatting from the insertion point, which we don't want here).

$d = $word.ActiveDocument

$contactInfo = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

$d.Content.Find.Execute("Dheeraj Chand", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dheeraj Chand^p" + $contactInfo, 2)

Write-Output "Inserted contact info paragraph after name heading."
